$wb = $excel.ActiveWorkbook

$backlog = $wb.Worksheets.Item("Backlog")
$backlog.Range("A27").Value = 2
$backlog.Range("E27").Value = "Planned"
$backlog.Range("A30").Value = 2
$backlog.Range("E30").Value = "Planned"

$sprint2 = $wb.Worksheets.Item("Sprint2")
$sprint2.Range("A4").Value = "US09"
$sprint2.Range("B4").Value = "Birth before death of parents"
$sprint2.Range("C4").Value = "DN"
$sprint2.Range("D4").Value = "Planned"
$sprint2.Range("E4").Value = 60
$sprint2.Range("F4").Value = 90

$sprint2.Range("A5").Value = "US22"
$sprint2.Range("B5").Value = "Unique IDs"
$sprint2.Range("C5").Value = "DN"
$sprint2.Range("D5").Value = "Planned"
$sprint2.Range("E5").Value = 40
$sprint2.Range("F5").Value = 90
